# Updated glasgow tribunal office
# The Glasgow Employment Tribunal moved premises: "Eagle Building, 215 Bothwell
# Street, Glasgow, G2 7TS" -> "Glasgow Tribunals Centre, 3 Atlantic Quay,
# 20 York Street, Glasgow, G2 8GT" (new 3rd address line added).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the extra "tribunalGlasgowAddressLine3" field, just
# above the existing "tribunalGlasgowTown" row (row 14, pushing everything
# below down by one).
$ws.Rows.Item(14).Insert()

# --- Glasgow address block (rows 12-17 after the insert) -----------------
# AddressLine1
$ws.Cells.Item(12, 2).Value = "Glasgow Tribunals Centre"
$ws.Cells.Item(12, 2).Style = "Normal"

# AddressLine2
$ws.Cells.Item(13, 2).Value = "3 Atlantic Quay"
$ws.Cells.Item(13, 2).Style = "Normal"

# AddressLine3 (brand new row/field)
$ws.Cells.Item(14, 1).Value = "tribunalGlasgowAddressLine3"
$ws.Cells.Item(14, 2).Value = "20 York Street"
$ws.Cells.Item(14, 2).Style = "Normal"

# Town (row 15) is unchanged - still "Glasgow".

# PostCode (row 16)
$ws.Cells.Item(16, 2).Value = "G2 8GT"
$ws.Cells.Item(16, 2).Style = "Normal"

# Telephone (row 17) is unchanged.

# --- Fix up the mailto hyperlinks, which don't auto-shift with the insert --
$ws.Hyperlinks.Delete()

$mailLinks = @(
  @(11,  "mailto:Manchesteret@justice.gov.uk"),
  @(20,  "mailto:glasgowet@justice.gov.uk"),
  @(29,  "mailto:aberdeenet@justice.gov.uk"),
  @(38,  "mailto:dundeeet@justice.gov.uk"),
  @(45,  "mailto:edinburghet@justice.gov.uk"),
  @(53,  "mailto:bristolet@justice.gov.uk"),
  @(61,  "mailto:LeedsET@justice.gov.uk"),
  @(70,  "mailto:londoncentralet@hmcts.gsi.gov.uk"),
  @(78,  "mailto:eastlondon@justice.gov.uk"),
  @(87,  "mailto:londonsouthet@hmcts.gsi.gov.uk"),
  @(94,  "mailto:e.midlandseastet@justice.gov.uk"),
  @(102, "mailto:MidlandsWestET@justice.gov.uk"),
  @(111, "mailto:newcastleet@hmcts.gsi.gov.uk"),
  @(118, "mailto:cardiffet@justice.gov.uk"),
  @(127, "mailto:watfordet@justice.gov.uk")
)
foreach ($link in $mailLinks) {
    $ws.Hyperlinks.Add($ws.Cells.Item($link[0], 2), $link[1]) | Out-Null
}

# --- Minor view/formatting touch-ups matching the resave ------------------
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30

$ws.Range("C16").Select()
